$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.184.99'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '1.826.38'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6038'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07111'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.07%  '
$ws.Range('E9').Value = '  -2.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07640'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').Value = '1.827.79'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.780'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6398'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009816'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '79.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').Value = '2.053.84'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.947'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.72%  '
$ws.Range('D19').Value = '29.187.39'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '230.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.040'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9987'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.035'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1278'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06751'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.450'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.456'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.805'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.773'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.134'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.717'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6588'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.530'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').Value = '1.234.04'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.756'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('E40').Value = '  -5.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.528'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9261'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '1.971.11'
$ws.Range('E44').Value = '  -2.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.87%  '
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.626'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.564'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05578'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.466'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.62%  '
